$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: update title (D23) and link (E23)
$ws.Range("D23").Value = "ML 고수분들께 질문드립니다!`n딥러닝 공부를 해오면서 요즘 더욱 더 기본기의 중요성을 느끼고 있는데 기본기를 직접 구현을 통해 복습해보려고 하는"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2698"

# Row 39: update title (D39) and link (E39)
$ws.Range("D39").Value = "Face Recognition with Dlib in Python"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Face-Recognition-with-Dlib-in-Python-1"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "아인트호벤 삼각형 (Einthoven’s Triangle)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/373"
